# Slide 1 ("The power of plain text") currently only has a Title and a
# Subtitle shape. The published deck also shows a Date placeholder
# (inherited from the Title Slide layout, idx=10, type="dt", sz="half")
# containing the literal text "2024-04-11". Turn that layout placeholder
# on for this slide and fill it in.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ppPlaceholderDate = 16 on the Title Slide layout - switching the
# slide's date footer element on materializes the corresponding
# "Date Placeholder" shape on the slide itself.
$dt = $s.HeadersFooters.DateAndTime
$dt.Visible = -1
$dt.UseFormat = 0

# Find the shape that was just materialized (the date placeholder) -
# look it up by placeholder type rather than assuming a fixed index.
$dateShape = $null
$phs = $s.Shapes.Placeholders
for ($i = 1; $i -le $phs.Count; $i++) {
    $candidate = $phs.Item($i)
    if ($candidate.PlaceholderFormat.Type -eq 16) {
        $dateShape = $candidate
    }
}

$dateShape.Name = "Date Placeholder 3"

$tr = $dateShape.TextFrame.TextRange
$tr.Text = "2024-04-11"

# Match the paragraph formatting used everywhere else in this deck
# (no bullet, explicit 0 margin/indent at level 1).
$tr.IndentLevel = 1
$tr.ParagraphFormat.Bullet.Visible = 0
$lvl = $dateShape.TextFrame.Ruler.Levels(1)
$lvl.LeftMargin = 0
$lvl.FirstMargin = 0
